# Apply weekly update: insert 3 new rows (Primera/Segunda/Tercera) for a new
# reporting date (2021-11-11, serial 44511) right before the existing block
# of Betarraga - Lo Valledor records, pushing the rest of that block (and the
# sheet) down by three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 455; this shifts the previous rows 455:565 down
# to 458:568, which is exactly the shift observed between the old and new
# data (row N in the new sheet == row N-3 in the old sheet, for N >= 458).
$ws.Rows("455:457").Insert()

# Common (constant) field values for these three new records.
$mercadoId   = 6
$mercado     = "Mercado Mayorista Lo Valledor de Santiago"
$region      = "Metropolitana"
$fecha       = 44511
$codreg      = 13
$categoriaId = 100114014
$categoria   = "Betarraga"
$variedad    = "Sin especificar"
$unidad      = "$/unidad"
$origen      = "Región Metropolitana"
$kgOUnidades = 1
$clasif      = "Hortaliza"

# Row 455 - Calidad "Primera"
$ws.Range("A455").Value = $mercadoId
$ws.Range("B455").Value = $mercado
$ws.Range("C455").Value = $region
$ws.Range("D455").Value = $fecha
$ws.Range("E455").Value = $codreg
$ws.Range("F455").Value = $categoriaId
$ws.Range("G455").Value = $categoria
$ws.Range("H455").Value = $variedad
$ws.Range("I455").Value = "Primera"
$ws.Range("J455").Value = 50000
$ws.Range("K455").Value = 90
$ws.Range("L455").Value = 100
$ws.Range("M455").Value = 94
$ws.Range("N455").Value = $unidad
$ws.Range("O455").Value = $origen
$ws.Range("P455").Value = 94
$ws.Range("Q455").Value = $kgOUnidades
$ws.Range("R455").Value = $clasif

# Row 456 - Calidad "Segunda"
$ws.Range("A456").Value = $mercadoId
$ws.Range("B456").Value = $mercado
$ws.Range("C456").Value = $region
$ws.Range("D456").Value = $fecha
$ws.Range("E456").Value = $codreg
$ws.Range("F456").Value = $categoriaId
$ws.Range("G456").Value = $categoria
$ws.Range("H456").Value = $variedad
$ws.Range("I456").Value = "Segunda"
$ws.Range("J456").Value = 43000
$ws.Range("K456").Value = 75
$ws.Range("L456").Value = 85
$ws.Range("M456").Value = 79
$ws.Range("N456").Value = $unidad
$ws.Range("O456").Value = $origen
$ws.Range("P456").Value = 79
$ws.Range("Q456").Value = $kgOUnidades
$ws.Range("R456").Value = $clasif

# Row 457 - Calidad "Tercera"
$ws.Range("A457").Value = $mercadoId
$ws.Range("B457").Value = $mercado
$ws.Range("C457").Value = $region
$ws.Range("D457").Value = $fecha
$ws.Range("E457").Value = $codreg
$ws.Range("F457").Value = $categoriaId
$ws.Range("G457").Value = $categoria
$ws.Range("H457").Value = $variedad
$ws.Range("I457").Value = "Tercera"
$ws.Range("J457").Value = 17000
$ws.Range("K457").Value = 60
$ws.Range("L457").Value = 60
$ws.Range("M457").Value = 60
$ws.Range("N457").Value = $unidad
$ws.Range("O457").Value = $origen
$ws.Range("P457").Value = 60
$ws.Range("Q457").Value = $kgOUnidades
$ws.Range("R457").Value = $clasif
